$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 20
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AN2").Value = 9.199999999999999
$ws.Range("F2").Value = 1.96
$ws.Range("O2").Value = 1.15
$ws.Range("P2").Value = 2.88
$ws.Range("T2").Value = 1.48
$ws.Range("U2").Value = 2.86
$ws.Range("F3").Value = 1.01
$ws.Range("H4").Value = 1.61
$ws.Range("F5").Value = 2.04
$ws.Range("H5").Value = 3.7
$ws.Range("J5").Value = 1.09
$ws.Range("K5").Value = 4.2
$ws.Range("N5").Value = 2.04
$ws.Range("O5").Value = 1.33
$ws.Range("Q5").Value = 2.6
$ws.Range("S5").Value = 2.6
$ws.Range("V5").Value = 1.2
$ws.Range("F6").Value = 1.33
$ws.Range("G6").Value = 1.68
$ws.Range("H6").Value = 7
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 980
$ws.Range("N6").Value = 2.44
$ws.Range("P6").Value = 1.54
$ws.Range("R6").Value = 1.16
$ws.Range("W6").Value = 2.46
$ws.Range("AN7").Value = 2.76
$ws.Range("H7").Value = 21
$ws.Range("P7").Value = 2.36
$ws.Range("U7").Value = 1.45
$ws.Range("W7").Value = 7.2
$ws.Range("P8").Value = 2.52
$ws.Range("S8").Value = 1.82
$ws.Range("AN9").Value = 980
$ws.Range("N9").Value = 1.94
$ws.Range("T9").Value = 1.71
$ws.Range("U9").Value = 2.14
$ws.Range("F10").Value = 2.34
$ws.Range("H10").Value = 2.96
$ws.Range("K10").Value = 3.95
$ws.Range("Q10").Value = 1.73
$ws.Range("F11").Value = 1.98
$ws.Range("G11").Value = 2.2
$ws.Range("U11").Value = 2.66
$ws.Range("F12").Value = 2.28
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 2.82
$ws.Range("O12").Value = 1.32
$ws.Range("S12").Value = 3.45
$ws.Range("V12").Value = 1.41
$ws.Range("X12").Value = 990
$ws.Range("F13").Value = 6.8
$ws.Range("I13").Value = 1.5
$ws.Range("K13").Value = 5.7
$ws.Range("L13").Value = 1.21
$ws.Range("N13").Value = 2.44
$ws.Range("P13").Value = 2.44
$ws.Range("S13").Value = 1.97
$ws.Range("U13").Value = 1.04
$ws.Range("V13").Value = 3
$ws.Range("F14").Value = 2.14
$ws.Range("AN15").Value = 4.6
$ws.Range("F15").Value = 1.25
$ws.Range("H15").Value = 9.4
$ws.Range("I15").Value = 19.5
$ws.Range("J15").Value = 6.4
$ws.Range("K15").Value = 9.6
$ws.Range("R15").Value = 1.75
$ws.Range("U15").Value = 1.74
$ws.Range("AK16").Value = 980
$ws.Range("G16").Value = 2.26
$ws.Range("I16").Value = 4.5
$ws.Range("L16").Value = 1.3
$ws.Range("W16").Value = 1.79
$ws.Range("I17").Value = 1.77
$ws.Range("AB18").Value = 12.5
$ws.Range("AK18").Value = 34
$ws.Range("Q18").Value = 2.02
$ws.Range("T18").Value = 1.78
$ws.Range("X18").Value = 15
$ws.Range("H19").Value = 1.74
$ws.Range("L20").Value = 1.18
$ws.Range("U20").Value = 2.14
$ws.Range("AG21").Value = 25
$ws.Range("Q21").Value = 1.86
$ws.Range("U21").Value = 1.96
$ws.Range("Z21").Value = 12.5
$ws.Range("G22").Value = 2.46
$ws.Range("H22").Value = 3.6
$ws.Range("AA24").Value = 1000
$ws.Range("AL24").Value = 30
$ws.Range("AN24").Value = 3
$ws.Range("J24").Value = 7.8
$ws.Range("L24").Value = 1.13
$ws.Range("Q24").Value = 1.31
$ws.Range("T24").Value = 1.77
$ws.Range("U24").Value = 2.08
$ws.Range("AG25").Value = 12
$ws.Range("AO25").Value = 17.5
$ws.Range("F25").Value = 2.64
$ws.Range("G25").Value = 2.66
$ws.Range("H25").Value = 2.78
$ws.Range("I25").Value = 2.8
$ws.Range("L25").Value = 1.3
$ws.Range("P25").Value = 2.4
$ws.Range("Q25").Value = 1.69
$ws.Range("R25").Value = 1.57
$ws.Range("S25").Value = 2.7
$ws.Range("T25").Value = 1.57
$ws.Range("V25").Value = 1.55
$ws.Range("Y25").Value = 15
